$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (row 29) with a new "test" entry that mirrors the
# structure of the existing log rows (20-28): a testmail subject/body
# pair plus the processed-email metadata columns.

$newRow = 29

$ws.Range("A$newRow").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Range("B$newRow").Value = "Beste klant,`nBedankt voor uw e-mail. Het spijt me, maar het lijkt erop dat u per ongeluk een testmail heeft gestuurd. Als u daadwerkelijk 200 stuks M8-bouten RVS wilt bestellen voor Van Dijk, raad ik u aan om contact op te nemen met onze verkoopafdeling via [verkoop@email.com] of telefonisch via [telefoonnummer]. Zij helpen u graag verder met uw bestelling.`nMet vriendelijke groet,`n[Naam] Nederlandse e-mailassistent van <bedrijfsnaam>"
$ws.Range("C$newRow").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$ws.Range("D$newRow").Value = "mailmind.test@zohomail.eu"
$ws.Range("E$newRow").Value = "Bestelling / Levering"
$ws.Range("F$newRow").Value = "2025-07-29 22:06:31"
$ws.Range("G$newRow").Value = "Ja"
$ws.Range("H$newRow").Value = "Nee"
$ws.Range("I$newRow").Value = "Ja"
$ws.Range("J$newRow").Value = "Nee"
